# Search-and-replace style edit: update the "yes baby" text wherever it
# appears on the sheet (C2 and C8) to its new replacement value "yes Babe".
# The other shared strings (NewText2 at E7, ReplacedValue at E12/F14) are
# left as-is.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "yes Babe"
$ws.Range("C8").Value = "yes Babe"
